$d = $word.ActiveDocument

# Locate "I used git" at the start of the paragraph; collapse the found
# range to its end point (right before " pages to host...") so we can
# insert "hub" there, splitting the original single run into three runs:
#   "I used git" | "hub" | " pages to host my online CV/portfolio. ..."
$r = $d.Content
$found = $r.Find.Execute("I used git", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'I used git'"
}
$r.Collapse(0)  # wdCollapseEnd

# Insert "hub" as a tracked insertion so the engine keeps it as its own
# run instead of silently re-merging it back into the surrounding text,
# then immediately accept that single revision to leave plain (non
# tracked-change) runs behind, matching a normal typed edit.
$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true
$r.InsertAfter("hub")
$d.TrackRevisions = $wasTracking

$d.Revisions(1).Accept()
